$wb = $excel.ActiveWorkbook

# --- sheet "sim": errorPropTestEnable flag flips from 1 to 0 ---
$sim = $wb.Worksheets.Item("sim")
$sim.Range("B7").Value = 0

# --- sheet "general": add three new camera-offset-of-body parameters ---
$general = $wb.Worksheets.Item("general")

$general.Range("A42").Value = "rcbx"
$general.Range("B42").Value = 0
$general.Range("C42").Value = "m"
$general.Range("D42").Value = "X component of camera offset of body"
$general.Range("E42").Formula = "=B42"

$general.Range("A43").Value = "rcby"
$general.Range("B43").Value = 0
$general.Range("C43").Value = "m"
$general.Range("D43").Value = "Y component of camera offset of body"
$general.Range("E43").Formula = "=B43"

$general.Range("A44").Value = "rcbz"
$general.Range("B44").Value = 0
$general.Range("C44").Value = "m"
$general.Range("D44").Value = "Z component of camera offset of body"
$general.Range("E44").Formula = "=B44"
